$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.3
$ws.Range("B3").Value = 0.2
$ws.Range("B4").Value = 0.15
$ws.Range("B6").Value = 0.1
$ws.Range("B7").Value = 0.05
$ws.Range("B8").Value = 0.1
